# lpagg/examples/VDI_4655_houses_example.xlsx
# "Make code more independent of requested energy types"
#
# - EFH/MFH now get 'calculate' for Q_TWW_a / W_a (instead of being left
#   blank and silently auto-filled elsewhere).
# - "Buero" renamed to "Büro".
# - house_type codes "G1G"/"G4G" renamed to "GHD/G1"/"GHA/G4".
# - copies/sigma values for Handel bumped (24->25, 4->5, 9->10).
# - Active cell/selection moved to D2.
# - PageSetup gets an explicit paper size / orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")

# --- Header row: rename "Buero" -> "Büro" ---
$ws.Range("D1").Value = "Büro"

# --- Row 4 (Q_TWW_a): EFH & MFH defaults now computed via 'calculate' ---
$ws.Range("B4").Value = "calculate"
$ws.Range("C4").Value = "calculate"

# --- Row 5 (W_a): EFH & MFH defaults now computed via 'calculate' ---
$ws.Range("B5").Value = "calculate"
$ws.Range("C5").Value = "calculate"

# --- Row 7 (house_type): renamed commercial/industrial VDI profile codes ---
$ws.Range("D7").Value = "GHD/G1"
$ws.Range("E7").Value = "GHA/G4"

# --- Row 10 (copies): Handel counts increased ---
$ws.Range("B10").Value = 25
$ws.Range("C10").Value = 25
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 10

# --- Selection moved from F15 to D2 ---
$ws.Range("D2").Select()

# --- Explicit page setup (paper size 9 = A4, portrait orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
